# Updated fastq files and related metadata
# - add 4 new data rows (38-41) to Sheet1
# - strip the ad-hoc header/column formatting (bold header row, bordered
#   cells, G/H column styling) back to the workbook default style
# - update the active-cell selection / used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Append the new rows of rnaSample metadata.
#    Columns: A harvestDate, B harvester, C bioSampleNumber, D rnaDate,
#             E rnaPreparer, F rnaSampleNumber, G rnaPrepMethod,
#             H roboticRNAPrep
#    A leading "'" forces these date-/bool-looking strings to stay text
#    (matching the existing shared-string cells) instead of being
#    auto-converted to a date serial or a boolean by COM's type sniffing.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 38; A = "08.09.18"; B = "H.BROWN"; C = 38; D = "08.13.18"; E = "H.BROWN"; F = 38; G = "TRIzol"; H = "False" },
    @{ Row = 39; A = "10.15.18"; B = "H.BROWN"; C = 39; D = "10.16.18"; E = "H.BROWN"; F = 39; G = "TRIzol"; H = "False" },
    @{ Row = 40; A = "10.30.18"; B = "H.BROWN"; C = 40; D = "11.01.18"; E = "H.BROWN"; F = 40; G = "TRIzol"; H = "False" },
    @{ Row = 41; A = "10.15.18"; B = "H.BROWN"; C = 41; D = "10.16.18"; E = "H.BROWN"; F = 41; G = "TRIzol"; H = "False" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = "'" + $r.A
    $ws.Range("B$row").Value = "'" + $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = "'" + $r.D
    $ws.Range("E$row").Value = "'" + $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = "'" + $r.G
    $ws.Range("H$row").Value = "'" + $r.H
    # Drop the quote-prefix / any inherited formatting so the new rows
    # end up on the plain default style, same as the rest of the sheet.
    $ws.Range("A$row" + ":H$row").ClearFormats()
}

# ---------------------------------------------------------------------
# 2) Remove the bespoke styling that used to live on the header row and
#    on the G/H data columns (bold+border+center on row 1, a plain font
#    override on G2:G37, a text-numfmt override on H2:H38). Clearing
#    formats drops the cell's style index back to the workbook default.
# ---------------------------------------------------------------------
$ws.Range("A1:H1").ClearFormats()
$ws.Range("G2:G37").ClearFormats()
$ws.Range("H2:H38").ClearFormats()

# ---------------------------------------------------------------------
# 3) Update the view: selection moves to I36, used range grows to H41.
# ---------------------------------------------------------------------
$ws.Range("I36").Select()

# Best-effort window resize/reposition (mirrors the workbookView change
# in workbook.xml: xWindow/yWindow -> 0/0, windowWidth/Height -> 33600/21000).
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 0
$win.Width = 33600
$win.Height = 21000
